# Grado08.xlsx - "Seguimiento" sheet update
# Commit: Fechas en seguimiento y corrección CN_04_08_CO
# Adds tracking dates to rows 6, 9 and 10, clears the two stray status
# comments that had been typed into column F for rows 9/10 (their text
# now lives only in rows 11/13, which keep their original wording), and
# moves the active selection to E12:E13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 - first tracking entry gets its publishing/manuscript dates
$ws.Range("B6").Value = 42092
$ws.Range("C6").Value = 42092

# Row 9 - full date trail + clear the leftover free-text note in F
$ws.Range("B9").Value = 42061
$ws.Range("C9").Value = 42068
$ws.Range("D9").Value = 42072
$ws.Range("E9").Value = 42072
$ws.Range("F9").Value = ""
$ws.Rows.Item(9).RowHeight = 16.5

# Row 10 - full date trail + clear the leftover free-text note in F
$ws.Range("B10").Value = 42083
$ws.Range("C10").Value = 42083
$ws.Range("D10").Value = 42101
$ws.Range("E10").Value = 42102
$ws.Range("F10").Value = ""

# Move the active selection to E12:E13
$ws.Range("E12:E13").Select() | Out-Null
